$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.4898806666666666
$ws.Range("H2").Value = 1.469642
$ws.Range("I2").Value = 0.01965582386814743
$ws.Range("J2").Value = 0.02319449744266509
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.500998666666666
$ws.Range("N2").Value = 7.502996
$ws.Range("O2").Value = 0.08179813614512804
$ws.Range("P2").Value = 0.1013554923061644
$ws.Range("Q2").Value = 1.225190894159111
$ws.Range("R2").Value = 11.026718047432
$ws.Range("S2").Value = 0.001607809756811381
$ws.Range("T2").Value = 0.002350889707095391
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.4898806666666666
$ws.Range("H3").Value = 1.469642
$ws.Range("I3").Value = 0.01965582386814743
$ws.Range("J3").Value = 0.02319449744266509
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.8623146666666667
$ws.Range("N3").Value = 2.586944
$ws.Range("O3").Value = 0.02820302683245761
$ws.Range("P3").Value = 0.03494617119461054
$ws.Range("Q3").Value = 0.4224312837831111
$ws.Range("R3").Value = 3.801881554048
$ws.Range("S3").Value = 0.0005543537279674227
$ws.Range("T3").Value = 0.0008105588784043308
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.4898806666666666
$ws.Range("H4").Value = 1.469642
$ws.Range("I4").Value = 0.01965582386814743
$ws.Range("J4").Value = 0.02319449744266509
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.273522333333333
$ws.Range("N4").Value = 12.820567
$ws.Range("O4").Value = 0.1397706309484552
$ws.Range("P4").Value = 0.1731888008375808
$ws.Range("Q4").Value = 2.093515969668222
$ws.Range("R4").Value = 18.841643727014
$ws.Range("S4").Value = 0.002747306903862672
$ws.Range("T4").Value = 0.004017027198125502
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.4898806666666666
$ws.Range("H5").Value = 1.469642
$ws.Range("I5").Value = 0.01965582386814743
$ws.Range("J5").Value = 0.02319449744266509
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.239195333333333
$ws.Range("N5").Value = 15.717586
$ws.Range("O5").Value = 0.1713541150096253
$ws.Range("P5").Value = 0.2123236726894799
$ws.Range("Q5").Value = 2.566580502690222
$ws.Range("R5").Value = 23.099224524212
$ws.Range("S5").Value = 0.003368106303711473
$ws.Range("T5").Value = 0.004924740883213402
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.4898806666666666
$ws.Range("H6").Value = 1.469642
$ws.Range("I6").Value = 0.01965582386814743
$ws.Range("J6").Value = 0.02319449744266509
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 17.6992215
$ws.Range("N6").Value = 35.398443
$ws.Range("O6").Value = 0.5788740910643339
$ws.Range("P6").Value = 0.4781858629721645
$ws.Range("Q6").Value = 8.670506427901
$ws.Range("R6").Value = 52.023038567406
$ws.Range("S6").Value = 0.01137824717579448
$ws.Range("T6").Value = 0.01109128077582647
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 13.025931
$ws.Range("H7").Value = 39.077793
$ws.Range("I7").Value = 0.5226485200912362
$ws.Range("J7").Value = 0.6167418798615554
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.500998666666666
$ws.Range("N7").Value = 7.502996
$ws.Range("O7").Value = 0.08179813614512804
$ws.Range("P7").Value = 0.1013554923061644
$ws.Range("Q7").Value = 32.57783606309199
$ws.Range("R7").Value = 293.200524567828
$ws.Range("S7").Value = 0.04275167480247263
$ws.Range("T7").Value = 0.06251017685919723
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 13.025931
$ws.Range("H8").Value = 39.077793
$ws.Range("I8").Value = 0.5226485200912362
$ws.Range("J8").Value = 0.6167418798615554
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.8623146666666667
$ws.Range("N8").Value = 2.586944
$ws.Range("O8").Value = 0.02820302683245761
$ws.Range("P8").Value = 0.03494617119461054
$ws.Range("Q8").Value = 11.232451348288
$ws.Range("R8").Value = 101.092062134592
$ws.Range("S8").Value = 0.0147402702360774
$ws.Range("T8").Value = 0.02155276731652784
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 13.025931
$ws.Range("H9").Value = 39.077793
$ws.Range("I9").Value = 0.5226485200912362
$ws.Range("J9").Value = 0.6167418798615554
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.273522333333333
$ws.Range("N9").Value = 12.820567
$ws.Range("O9").Value = 0.1397706309484552
$ws.Range("P9").Value = 0.1731888008375808
$ws.Range("Q9").Value = 55.666607040959
$ws.Range("R9").Value = 500.999463368631
$ws.Range("S9").Value = 0.07305091341742846
$ws.Range("T9").Value = 0.1068127865995381
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 13.025931
$ws.Range("H10").Value = 39.077793
$ws.Range("I10").Value = 0.5226485200912362
$ws.Range("J10").Value = 0.6167418798615554
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.239195333333333
$ws.Range("N10").Value = 15.717586
$ws.Range("O10").Value = 0.1713541150096253
$ws.Range("P10").Value = 0.2123236726894799
$ws.Range("Q10").Value = 68.24539690752199
$ws.Range("R10").Value = 614.208572167698
$ws.Range("S10").Value = 0.08955797462132414
$ws.Range("T10").Value = 0.1309489010336194
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 13.025931
$ws.Range("H11").Value = 39.077793
$ws.Range("I11").Value = 0.5226485200912362
$ws.Range("J11").Value = 0.6167418798615554
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 17.6992215
$ws.Range("N11").Value = 35.398443
$ws.Range("O11").Value = 0.5788740910643339
$ws.Range("P11").Value = 0.4781858629721645
$ws.Range("Q11").Value = 230.5488380127165
$ws.Range("R11").Value = 1383.293028076299
$ws.Range("S11").Value = 0.3025476870139336
$ws.Range("T11").Value = 0.2949172480526729
$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 11.4071155
$ws.Range("H12").Value = 22.814231
$ws.Range("I12").Value = 0.4576956560406163
$ws.Range("J12").Value = 0.3600636226957795
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 2.500998666666666
$ws.Range("N12").Value = 7.502996
$ws.Range("O12").Value = 0.08179813614512804
$ws.Range("P12").Value = 0.1013554923061644
$ws.Range("Q12").Value = 28.52918065601266
$ws.Range("R12").Value = 171.175083936076
$ws.Range("S12").Value = 0.03743865158584403
$ws.Range("T12").Value = 0.03649442573987175
$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 11.4071155
$ws.Range("H13").Value = 22.814231
$ws.Range("I13").Value = 0.4576956560406163
$ws.Range("J13").Value = 0.3600636226957795
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.8623146666666667
$ws.Range("N13").Value = 2.586944
$ws.Range("O13").Value = 0.02820302683245761
$ws.Range("P13").Value = 0.03494617119461054
$ws.Range("Q13").Value = 9.836523000010667
$ws.Range("R13").Value = 59.019138000064
$ws.Range("S13").Value = 0.01290840286841279
$ws.Range("T13").Value = 0.01258284499967837
$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 11.4071155
$ws.Range("H14").Value = 22.814231
$ws.Range("I14").Value = 0.4576956560406163
$ws.Range("J14").Value = 0.3600636226957795
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 4.273522333333333
$ws.Range("N14").Value = 12.820567
$ws.Range("O14").Value = 0.1397706309484552
$ws.Range("P14").Value = 0.1731888008375808
$ws.Range("Q14").Value = 48.74856284816283
$ws.Range("R14").Value = 292.491377088977
$ws.Range("S14").Value = 0.06397241062716409
$ws.Range("T14").Value = 0.06235898703991719
$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 11.4071155
$ws.Range("H15").Value = 22.814231
$ws.Range("I15").Value = 0.4576956560406163
$ws.Range("J15").Value = 0.3600636226957795
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 5.239195333333333
$ws.Range("N15").Value = 15.717586
$ws.Range("O15").Value = 0.1713541150096253
$ws.Range("P15").Value = 0.2123236726894799
$ws.Range("Q15").Value = 59.76410629439433
$ws.Range("R15").Value = 358.584637766366
$ws.Range("S15").Value = 0.07842803408458968
$ws.Range("T15").Value = 0.07645003077264707
$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 11.4071155
$ws.Range("H16").Value = 22.814231
$ws.Range("I16").Value = 0.4576956560406163
$ws.Range("J16").Value = 0.3600636226957795
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 17.6992215
$ws.Range("N16").Value = 35.398443
$ws.Range("O16").Value = 0.5788740910643339
$ws.Range("P16").Value = 0.4781858629721645
$ws.Range("Q16").Value = 201.8970639105833
$ws.Range("R16").Value = 807.588255642333
$ws.Range("S16").Value = 0.2649481568746058
$ws.Range("T16").Value = 0.1721773341436652
